$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 57:58, shifting existing rows 57-129 down to 59-131
$ws.Rows("57:58").Insert()

# Row 57: new weekly record (Primera)
$ws.Range("A57").Value = 2
$ws.Range("B57").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C57").Value = "Coquimbo"
$ws.Range("D57").Value = 44601
$ws.Range("E57").Value = 4
$ws.Range("F57").Value = 100112043
$ws.Range("G57").Value = "Pepino ensalada"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 600
$ws.Range("K57").Value = 11000
$ws.Range("L57").Value = 12000
$ws.Range("M57").Value = 11500
$ws.Range("N57").Value = "$/caja 70 unidades"
$ws.Range("O57").Value = "Provincia de Limarí"
$ws.Range("P57").Value = 164
$ws.Range("Q57").Value = 70
$ws.Range("R57").Value = "Hortaliza"

# Row 58: new weekly record (Segunda)
$ws.Range("A58").Value = 2
$ws.Range("B58").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C58").Value = "Coquimbo"
$ws.Range("D58").Value = 44601
$ws.Range("E58").Value = 4
$ws.Range("F58").Value = 100112043
$ws.Range("G58").Value = "Pepino ensalada"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Segunda"
$ws.Range("J58").Value = 500
$ws.Range("K58").Value = 9000
$ws.Range("L58").Value = 10000
$ws.Range("M58").Value = 9500
$ws.Range("N58").Value = "$/caja 100 unidades"
$ws.Range("O58").Value = "Provincia de Limarí"
$ws.Range("P58").Value = 95
$ws.Range("Q58").Value = 100
$ws.Range("R58").Value = "Hortaliza"
